# Updated cryptos list on Thu Nov 23 13:53:20 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to keep a literal text representation (avoid Excel
    # auto-converting numeric-looking strings like "0.615" or "1.00" into
    # real numbers, which would lose trailing zeros / formatting).
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# --- Row 7 / Row 8 swap: Solana <-> USDC ---
Set-TextValue "B7" "USDC"
Set-TextValue "C7" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.06%  "

Set-TextValue "B8" "Solana"
Set-TextValue "C8" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D8" "57.66"
Set-TextValue "E8" "  +4.41%  "

# --- Row 39 / Row 40 swap: WEMIXToken <-> RenderToken ---
Set-TextValue "B39" "RenderToken"
Set-TextValue "C39" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D39" "3.26"
Set-TextValue "E39" "  -1.82%  "

Set-TextValue "B40" "WEMIXToken"
Set-TextValue "C40" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D40" "1.74"
Set-TextValue "E40" "  -1.46%  "

# --- Row 45 / Row 46 swap: Cronos <-> TrustWalletToken ---
Set-TextValue "B45" "TrustWalletToken"
Set-TextValue "C45" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D45" "1.16"
Set-TextValue "E45" "  +4.22%  "

Set-TextValue "B46" "Cronos"
Set-TextValue "C46" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D46" "0.0923"
Set-TextValue "E46" "  -2.25%  "

# --- Plain value/percentage updates (Price column D, Volume(1h) column E) ---

# Row 2 - Bitcoin
Set-TextValue "D2" "37.282.38"
Set-TextValue "E2" "  +1.78%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.052.85"
Set-TextValue "E3" "  +0.89%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.19%  "

# Row 5 - BNB
Set-TextValue "D5" "232.37"
Set-TextValue "E5" "  -0.81%  "

# Row 6 - XRP
Set-TextValue "D6" "0.615"
Set-TextValue "E6" "  +2.38%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.382"
Set-TextValue "E9" "  +2.53%  "

# Row 10 - OKB
Set-TextValue "D10" "58.28"
Set-TextValue "E10" "  +1.23%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0757"
Set-TextValue "E11" "  +0.59%  "

# Row 12 - TRON
Set-TextValue "E12" "  +1.16%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "2.358.57"
Set-TextValue "E13" "  +1.23%  "

# Row 14 - Chainlink
Set-TextValue "D14" "14.35"
Set-TextValue "E14" "  +0.07%  "

# Row 15 - Avalanche
Set-TextValue "D15" "20.67"
Set-TextValue "E15" "  +1.48%  "

# Row 16 - Polygon
Set-TextValue "D16" "0.771"
Set-TextValue "E16" "  +0.92%  "

# Row 17 - Polkadot
Set-TextValue "D17" "5.15"
Set-TextValue "E17" "  +0.88%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.044.87"
Set-TextValue "E18" "  +0.47%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "37.232.75"
Set-TextValue "E19" "  +1.38%  "

# Row 20 - Uniswap
Set-TextValue "D20" "6.22"
Set-TextValue "E20" "  +13.58%  "

# Row 21 - Litecoin
Set-TextValue "D21" "69.14"
Set-TextValue "E21" "  +1.92%  "

# Row 22 - ShibaInu
Set-TextValue "D22" "0.0₃0809"
Set-TextValue "E22" "  +1.01%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "224.79"
Set-TextValue "E23" "  +1.57%  "

# Row 24 - Dai
Set-TextValue "D24" "0.999"
Set-TextValue "E24" "  -0.22%  "

# Row 25 - Toncoin
Set-TextValue "E25" "  +2.86%  "

# Row 26 - PancakeSwap
Set-TextValue "D26" "2.38"
Set-TextValue "E26" "  -0.88%  "

# Row 27 - Monero
Set-TextValue "D27" "165.16"
Set-TextValue "E27" "  +1.17%  "

# Row 28 - ImmutableX
Set-TextValue "E28" "  +7.29%  "

# Row 29 - Cosmos
Set-TextValue "D29" "8.77"
Set-TextValue "E29" "  +1.39%  "

# Row 30 - Kaspa
Set-TextValue "E30" "  -5.00%  "

# Row 31 - EthereumClassic
Set-TextValue "D31" "18.99"
Set-TextValue "E31" "  -0.13%  "

# Row 32 - Stellar
Set-TextValue "D32" "0.117"
Set-TextValue "E32" "  +0.36%  "

# Row 33 - Filecoin
Set-TextValue "D33" "4.47"
Set-TextValue "E33" "  +1.96%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.0614"
Set-TextValue "E34" "  +1.72%  "

# Row 35 - LidoDAOToken
Set-TextValue "D35" "2.53"
Set-TextValue "E35" "  +2.58%  "

# Row 36 - InternetComputer(DFINITY)
Set-TextValue "D36" "4.53"
Set-TextValue "E36" "  +6.00%  "

# Row 37 - BinanceUSD
Set-TextValue "E37" "  +0.18%  "

# Row 38 - THORChain
Set-TextValue "D38" "5.84"
Set-TextValue "E38" "  +0.78%  "

# Row 41 - FTXToken
Set-TextValue "D41" "4.71"
Set-TextValue "E41" "  +12.58%  "

# Row 42 - HuobiToken
Set-TextValue "E42" "  +0.99%  "

# Row 43 - Maker
Set-TextValue "D43" "1.475.06"
Set-TextValue "E43" "  +0.97%  "

# Row 44 - Aave
Set-TextValue "D44" "96.48"
Set-TextValue "E44" "  +6.02%  "

# Row 47 - VeChain
Set-TextValue "E47" "  +2.50%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "15.48"
Set-TextValue "E48" "  -1.02%  "

# Row 49 - ARBITRUM
Set-TextValue "E49" "  +0.57%  "

# Row 50 - FraxShare
Set-TextValue "D50" "7.14"
Set-TextValue "E50" "  +3.54%  "

# Row 51 - MXToken
Set-TextValue "E51" "  +1.75%  "
